$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Fill in the "unplaced" layout block (rows 1-17) on Sheet3 - the static
# object placements (lockers, medkits, offices, etc.) for each room.
$ws3.Range('L1').Value2 = 'A'
$ws3.Range('O1').Value2 = 'B'
$ws3.Range('D2').Value2 = 'end'
$ws3.Range('H2').Value2 = 'end'
$ws3.Range('J2').Value2 = 'end'
$ws3.Range('L2').Value2 = 'h'
$ws3.Range('O2').Value2 = 'h'
$ws3.Range('B3').Value2 = 'lr'
$ws3.Range('C3').Value2 = 'h'
$ws3.Range('D3').Value2 = 't'
$ws3.Range('E3').Value2 = 'h'
$ws3.Range('F3').Value2 = 'l'
$ws3.Range('G3').Value2 = 860
$ws3.Range('H3').Value2 = 4
$ws3.Range('I3').Value2 = 'h'
$ws3.Range('J3').Value2 = 't'
$ws3.Range('K3').Value2 = 'med'
$ws3.Range('L3').Value2 = 't'
$ws3.Range('M3').Value2 = 'h'
$ws3.Range('N3').Value2 = 'h'
$ws3.Range('O3').Value2 = 'elec'
$ws3.Range('B4').Value2 = 'h'
$ws3.Range('H4').Value2 = 'caf'
$ws3.Range('B5').Value2 = 'office'
$ws3.Range('H5').Value2 = 'h'
$ws3.Range('B6').Value2 = 'z'
$ws3.Range('E6').Value2 = 'z'
$ws3.Range('H6').Value2 = 'z'
$ws3.Range('L6').Value2 = 'z'
$ws3.Range('A7').Value2 = 8
$ws3.Range('B7').Value2 = 4
$ws3.Range('C7').Value2 = 'h'
$ws3.Range('D7').Value2 = 't'
$ws3.Range('E7').Value2 = 't'
$ws3.Range('F7').Value2 = 't'
$ws3.Range('G7').Value2 = 'h'
$ws3.Range('H7').Value2 = 't'
$ws3.Range('I7').Value2 = 'h'
$ws3.Range('J7').Value2 = 966
$ws3.Range('K7').Value2 = 'pump'
$ws3.Range('L7').Value2 = 4
$ws3.Range('M7').Value2 = 106
$ws3.Range('B8').Value2 = 'h'
$ws3.Range('D8').Value2 = 'h'
$ws3.Range('F8').Value2 = 'ltest'
$ws3.Range('J8').Value2 = 35
$ws3.Range('L8').Value2 = 895
$ws3.Range('B9').Value2 = 'h'
$ws3.Range('D9').Value2 = 'elev'
$ws3.Range('F9').Value2 = 'h'
$ws3.Range('O9').Value2 = 79
$ws3.Range('B10').Value2 = 't'
$ws3.Range('C10').Value2 = 'h'
$ws3.Range('D10').Value2 = 't'
$ws3.Range('E10').Value2 = 'bell'
$ws3.Range('F10').Value2 = 't'
$ws3.Range('G10').Value2 = 't'
$ws3.Range('H10').Value2 = 'h'
$ws3.Range('I10').Value2 = 'doc'
$ws3.Range('J10').Value2 = 't'
$ws3.Range('K10').Value2 = 'h'
$ws3.Range('L10').Value2 = 'h'
$ws3.Range('M10').Value2 = 'mt'
$ws3.Range('N10').Value2 = 'h'
$ws3.Range('O10').Value2 = 'c'
$ws3.Range('B11').Value2 = 'h'
$ws3.Range('E11').Value2 = 'nuke'
$ws3.Range('G11').Value2 = 'h'
$ws3.Range('J11').Value2 = 'h'
$ws3.Range('B12').Value2 = 'z'
$ws3.Range('E12').Value2 = 'z'
$ws3.Range('G12').Value2 = 'z'
$ws3.Range('J12').Value2 = 'z'
$ws3.Range('B13').Value2 = 't'
$ws3.Range('C13').Value2 = 'h'
$ws3.Range('D13').Value2 = 'h'
$ws3.Range('E13').Value2 = 't'
$ws3.Range('F13').Value2 = 'h'
$ws3.Range('G13').Value2 = 't'
$ws3.Range('H13').Value2 = 'h'
$ws3.Range('I13').Value2 = 1123
$ws3.Range('J13').Value2 = 't'
$ws3.Range('K13').Value2 = 1499
$ws3.Range('L13').Value2 = 'end'
$ws3.Range('B14').Value2 = 12
$ws3.Range('B15').Value2 = 'security'
$ws3.Range('D15').Value2 = 178
$ws3.Range('F15').Value2 = 205
$ws3.Range('H15').Value2 = 'arch'
$ws3.Range('J15').Value2 = 914
$ws3.Range('B16').Value2 = 1162
$ws3.Range('C16').Value2 = 'h'
$ws3.Range('D16').Value2 = 939
$ws3.Range('E16').Value2 = 'h'
$ws3.Range('F16').Value2 = 't'
$ws3.Range('G16').Value2 = 'h'
$ws3.Range('H16').Value2 = 4
$ws3.Range('I16').Value2 = 970
$ws3.Range('J16').Value2 = 't'
$ws3.Range('K16').Value2 = 'spcs'
$ws3.Range('L16').Value2 = 'h'
$ws3.Range('M16').Value2 = 'test'
$ws3.Range('N16').Value2 = 'closets'
$ws3.Range('O16').Value2 = 372
$ws3.Range('H17').Value2 = 173

# Duplicate the block 18 rows down (rows 19-35) to represent the "placed"
# copy of the same objects.
$ws3.Range("A1:O17").Copy()
$ws3.Range("A19").PasteSpecial(-4163)   # xlPasteAll

# Mark every non-blank cell in the duplicated block with the same black
# "done" fill already used elsewhere in the workbook (Sheet2!A2), without
# touching/adding any new style/fill definitions.
$rng = $ws3.Range("A19:O35").SpecialCells(2)   # xlCellTypeConstants
foreach ($area in $rng.Areas) {
  $ws2.Range("A2").Copy()
  $area.PasteSpecial(-4122)   # xlPasteFormats
}

# Sheet3 becomes the active sheet/tab, with Q17 selected.
$ws3.Activate()
$ws3.Range("Q17").Select()
